$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Briefly worked with the clang static analyzer"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "Briefly worked with the clang static analyzer, 1 year."
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null
